$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 17's timestamp (Fecha) value — small precision correction
$ws.Range("A17").Value = 45864.87532446759

# Append the new row 18 with the latest scheduled-task reading
$ws.Range("A18").Value = 45864.91690162032
$ws.Range("B18").Value = 2025
$ws.Range("C18").Value = 30
$ws.Range("D18").Value = 13.83
$ws.Range("E18").Value = 89.76000000000001
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 1.68
$ws.Range("H18").Value = "WSW"
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = "22:00:20"

# Match the date-formatted style used by the rest of column A
$ws.Range("A18").NumberFormat = $ws.Range("A17").NumberFormat
